$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new row at row 2 (shifts existing rows 2..21 down to 3..22)
$ws.Rows.Item(2).Insert()

# Excel's row insert copies the formatting of the row above (the bold header row);
# strip that back off so the new data row matches the other plain data rows.
$ws.Rows.Item(2).ClearFormats()

# 2) Populate the newly inserted row 2 with its data
$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = "walkingToRunning"
$ws.Cells.Item(2, 3).Value = 6.48236270504103
$ws.Cells.Item(2, 4).Value = -11.76500675559703
$ws.Cells.Item(2, 5).Value = 1.762938408561835
$ws.Cells.Item(2, 6).Value = 1.779657483100891
$ws.Cells.Item(2, 7).Value = 1.38677453994751
$ws.Cells.Item(2, 8).Value = -1.217397570610046

# 3) Append 10 new rows of data (rows 22..31) after the now-shifted last row (row 22 holds
#    the former last row's data, but with an updated timestamp of 2000 and continues the series)
$newRows = @(
    @(2000, "walkingToRunning", 7.328208127733031, -16.61346639454043, 21.1434863812355, 6.44285249710083, -5.062869071960449, 7.707920551300049),
    @(2100, "walkingToRunning", 13.14612888104349, -0.4819092987650375, 29.9319970146727, 0.4938832223415375, -9.879995346069336, -0.3310975134372711),
    @(2200, "walkingToRunning", 22.3981020806244, -18.41001325154172, 26.03648814838894, 6.862566947937012, 17.44670104980469, -6.492054462432861),
    @(2300, "walkingToRunning", 25.07217257589271, -46.9146286031819, 26.18006417501044, 0.584963321685791, -0.4454802870750427, 2.492385864257812),
    @(2400, "walkingToRunning", 1.007597878492922, 2.652457189823263, 3.160516032856393, -7.439141273498535, -3.388273239135742, 10.04511070251465),
    @(2500, "walkingToRunning", -11.37837042466059, -20.30985535047345, 3.219444780718657, 1.835583806037903, 12.65514183044434, -3.548196077346802),
    @(2600, "walkingToRunning", -19.19458270599844, -40.56597906581633, 23.86108243531293, -0.2848250865936279, -6.074337482452393, 3.493201971054077),
    @(2700, "walkingToRunning", -2.306829021780334, -12.78860808340835, 22.76510469426068, 7.885753154754639, -7.940680980682373, -2.813696384429932),
    @(2800, "walkingToRunning", -29.36129014663261, -62.85535855846575, 59.51861236108672, 1.459545493125916, 0.23096264898777, -2.083990097045898),
    @(2900, "walkingToRunning", -37.91703132502937, 6.169044086287685, -16.59176512449484, -3.29206657409668, 4.709334373474121, 3.162437200546265)
)

$startRow = 22
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $values = $newRows[$i]
    for ($c = 1; $c -le 8; $c++) {
        $ws.Cells.Item($r, $c).Value = $values[$c - 1]
    }
}
